$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 8 and 9: swap match details (F:V) ---
# Row 8
$ws.Cells.Item(8, 6).Value = "Dangkor"
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = "Prey Veng"
$ws.Cells.Item(8, 9).Value = 5
$ws.Cells.Item(8, 10).Value = 2.04
$ws.Cells.Item(8, 11).Value = "12/08/2023 05:12"
$ws.Cells.Item(8, 12).Value = 2.04
$ws.Cells.Item(8, 13).Value = "12/08/2023 12:07"
$ws.Cells.Item(8, 14).Value = 3.77
$ws.Cells.Item(8, 15).Value = "12/08/2023 05:12"
$ws.Cells.Item(8, 16).Value = 3.63
$ws.Cells.Item(8, 17).Value = "12/08/2023 12:07"
$ws.Cells.Item(8, 18).Value = 2.73
$ws.Cells.Item(8, 19).Value = "12/08/2023 05:12"
$ws.Cells.Item(8, 20).Value = 2.94
$ws.Cells.Item(8, 21).Value = "12/08/2023 12:07"
$ws.Cells.Item(8, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-prey-veng/b1Zt7ysa/"

# Row 9
$ws.Cells.Item(9, 6).Value = "NagaWorld"
$ws.Cells.Item(9, 7).Value = 1
$ws.Cells.Item(9, 8).Value = "Visakha"
$ws.Cells.Item(9, 9).Value = 4
$ws.Cells.Item(9, 10).Value = 2.75
$ws.Cells.Item(9, 11).Value = "11/08/2023 01:12"
$ws.Cells.Item(9, 12).Value = 2.52
$ws.Cells.Item(9, 13).Value = "12/08/2023 09:30"
$ws.Cells.Item(9, 14).Value = 3.53
$ws.Cells.Item(9, 15).Value = "11/08/2023 01:12"
$ws.Cells.Item(9, 16).Value = 3.64
$ws.Cells.Item(9, 17).Value = "12/08/2023 11:02"
$ws.Cells.Item(9, 18).Value = 2.03
$ws.Cells.Item(9, 19).Value = "11/08/2023 01:12"
$ws.Cells.Item(9, 20).Value = 2.28
$ws.Cells.Item(9, 21).Value = "12/08/2023 09:30"
$ws.Cells.Item(9, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/nagaworld-visakha/xhYp6ed5/"

# --- Update rows 38-41: rotate match details (F:V) ---
# Row 38
$ws.Cells.Item(38, 6).Value = "Visakha"
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = "Angkor Tiger"
$ws.Cells.Item(38, 9).Value = 1
$ws.Cells.Item(38, 10).Value = 1.15
$ws.Cells.Item(38, 11).Value = "21/10/2023 00:43"
$ws.Cells.Item(38, 12).Value = 1.19
$ws.Cells.Item(38, 13).Value = "22/10/2023 12:45"
$ws.Cells.Item(38, 14).Value = 6.29
$ws.Cells.Item(38, 15).Value = "21/10/2023 00:43"
$ws.Cells.Item(38, 16).Value = 6.51
$ws.Cells.Item(38, 17).Value = "22/10/2023 12:45"
$ws.Cells.Item(38, 18).Value = 8.19
$ws.Cells.Item(38, 19).Value = "21/10/2023 00:43"
$ws.Cells.Item(38, 20).Value = 8.7
$ws.Cells.Item(38, 21).Value = "22/10/2023 12:45"
$ws.Cells.Item(38, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/visakha-angkor-tiger/8YkmqVgG/"

# Row 39
$ws.Cells.Item(39, 6).Value = "Svay Rieng"
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(39, 8).Value = "Kirivong Sok Sen Chey"
$ws.Cells.Item(39, 9).Value = 1
$ws.Cells.Item(39, 10).Value = 1.11
$ws.Cells.Item(39, 11).Value = "21/10/2023 00:13"
$ws.Cells.Item(39, 12).Value = 1.25
$ws.Cells.Item(39, 13).Value = "22/10/2023 12:44"
$ws.Cells.Item(39, 14).Value = 7
$ws.Cells.Item(39, 15).Value = "21/10/2023 00:13"
$ws.Cells.Item(39, 16).Value = 5.69
$ws.Cells.Item(39, 17).Value = "22/10/2023 12:51"
$ws.Cells.Item(39, 18).Value = 9.71
$ws.Cells.Item(39, 19).Value = "21/10/2023 00:13"
$ws.Cells.Item(39, 20).Value = 7.79
$ws.Cells.Item(39, 21).Value = "22/10/2023 12:46"
$ws.Cells.Item(39, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/svay-rieng-kirivong-sok-sen-chey/vaoqpBvA/"

# Row 40
$ws.Cells.Item(40, 6).Value = "Dangkor"
$ws.Cells.Item(40, 7).Value = 2
$ws.Cells.Item(40, 8).Value = "Boeung Ket"
$ws.Cells.Item(40, 9).Value = 1
$ws.Cells.Item(40, 10).Value = 4.37
$ws.Cells.Item(40, 11).Value = "21/10/2023 00:13"
$ws.Cells.Item(40, 12).Value = 4.18
$ws.Cells.Item(40, 13).Value = "22/10/2023 12:44"
$ws.Cells.Item(40, 14).Value = 4.07
$ws.Cells.Item(40, 15).Value = "21/10/2023 00:13"
$ws.Cells.Item(40, 16).Value = 4.44
$ws.Cells.Item(40, 17).Value = "22/10/2023 12:44"
$ws.Cells.Item(40, 18).Value = 1.48
$ws.Cells.Item(40, 19).Value = "21/10/2023 00:13"
$ws.Cells.Item(40, 20).Value = 1.56
$ws.Cells.Item(40, 21).Value = "22/10/2023 12:44"
$ws.Cells.Item(40, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-boeung-ket/nPlirk9M/"

# Row 41
$ws.Cells.Item(41, 6).Value = "NagaWorld"
$ws.Cells.Item(41, 7).Value = 2
$ws.Cells.Item(41, 8).Value = "Tiffy Army"
$ws.Cells.Item(41, 9).Value = 1
$ws.Cells.Item(41, 10).Value = 1.84
$ws.Cells.Item(41, 11).Value = "21/10/2023 00:13"
$ws.Cells.Item(41, 12).Value = 1.95
$ws.Cells.Item(41, 13).Value = "22/10/2023 12:43"
$ws.Cells.Item(41, 14).Value = 3.45
$ws.Cells.Item(41, 15).Value = "21/10/2023 00:13"
$ws.Cells.Item(41, 16).Value = 3.73
$ws.Cells.Item(41, 17).Value = "22/10/2023 12:43"
$ws.Cells.Item(41, 18).Value = 3.27
$ws.Cells.Item(41, 19).Value = "21/10/2023 00:13"
$ws.Cells.Item(41, 20).Value = 3.07
$ws.Cells.Item(41, 21).Value = "22/10/2023 12:43"
$ws.Cells.Item(41, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/nagaworld-tiffy-army/69VNxTWq/"

# --- Add new rows 47, 48, 49 (copy formatting from row 46 first) ---
# Row 47
$ws.Range("A46:V46").Copy()
$ws.Range("A47:V47").PasteSpecial(-4122)
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "cambodia"
$ws.Cells.Item(47, 3).Value = "cpl"
$ws.Cells.Item(47, 4).Value = "2023-2024"
$ws.Cells.Item(47, 5).Value = 45234.40625
$ws.Cells.Item(47, 6).Value = "Kirivong Sok Sen Chey"
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = "Visakha"
$ws.Cells.Item(47, 9).Value = 2
$ws.Cells.Item(47, 10).Value = 5.22
$ws.Cells.Item(47, 11).Value = "02/11/2023 22:12"
$ws.Cells.Item(47, 12).Value = 5.47
$ws.Cells.Item(47, 13).Value = "04/11/2023 09:30"
$ws.Cells.Item(47, 14).Value = 4.49
$ws.Cells.Item(47, 15).Value = "02/11/2023 22:12"
$ws.Cells.Item(47, 16).Value = 4.73
$ws.Cells.Item(47, 17).Value = "04/11/2023 09:30"
$ws.Cells.Item(47, 18).Value = 1.36
$ws.Cells.Item(47, 19).Value = "02/11/2023 22:12"
$ws.Cells.Item(47, 20).Value = 1.4
$ws.Cells.Item(47, 21).Value = "04/11/2023 09:30"
$ws.Cells.Item(47, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/kirivong-sok-sen-chey-visakha/6TkmuMhr/"

# Row 48
$ws.Range("A46:V46").Copy()
$ws.Range("A48:V48").PasteSpecial(-4122)
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "cambodia"
$ws.Cells.Item(48, 3).Value = "cpl"
$ws.Cells.Item(48, 4).Value = "2023-2024"
$ws.Cells.Item(48, 5).Value = 45234.5
$ws.Cells.Item(48, 6).Value = "Boeung Ket"
$ws.Cells.Item(48, 7).Value = 5
$ws.Cells.Item(48, 8).Value = "Phnom Penh Crown"
$ws.Cells.Item(48, 9).Value = 5
$ws.Cells.Item(48, 10).Value = 3.47
$ws.Cells.Item(48, 11).Value = "03/11/2023 00:12"
$ws.Cells.Item(48, 12).Value = 4.21
$ws.Cells.Item(48, 13).Value = "04/11/2023 11:45"
$ws.Cells.Item(48, 14).Value = 3.65
$ws.Cells.Item(48, 15).Value = "03/11/2023 00:12"
$ws.Cells.Item(48, 16).Value = 3.9
$ws.Cells.Item(48, 17).Value = "04/11/2023 11:46"
$ws.Cells.Item(48, 18).Value = 1.7
$ws.Cells.Item(48, 19).Value = "03/11/2023 00:12"
$ws.Cells.Item(48, 20).Value = 1.63
$ws.Cells.Item(48, 21).Value = "04/11/2023 11:46"
$ws.Cells.Item(48, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/boeung-ket-phnom-penh-crown/Eggew0Nf/"

# Row 49
$ws.Range("A46:V46").Copy()
$ws.Range("A49:V49").PasteSpecial(-4122)
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "cambodia"
$ws.Cells.Item(49, 3).Value = "cpl"
$ws.Cells.Item(49, 4).Value = "2023-2024"
$ws.Cells.Item(49, 5).Value = 45234.5
$ws.Cells.Item(49, 6).Value = "NagaWorld"
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = "Prey Veng"
$ws.Cells.Item(49, 9).Value = 4
$ws.Cells.Item(49, 10).Value = 1.6
$ws.Cells.Item(49, 11).Value = "03/11/2023 00:12"
$ws.Cells.Item(49, 12).Value = 1.89
$ws.Cells.Item(49, 13).Value = "04/11/2023 11:55"
$ws.Cells.Item(49, 14).Value = 3.82
$ws.Cells.Item(49, 15).Value = "03/11/2023 00:12"
$ws.Cells.Item(49, 16).Value = 3.69
$ws.Cells.Item(49, 17).Value = "04/11/2023 11:55"
$ws.Cells.Item(49, 18).Value = 3.8
$ws.Cells.Item(49, 19).Value = "03/11/2023 00:12"
$ws.Cells.Item(49, 20).Value = 3.26
$ws.Cells.Item(49, 21).Value = "04/11/2023 11:55"
$ws.Cells.Item(49, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/nagaworld-prey-veng/QVyJdlOF/"

$excel.CutCopyMode = $false
